$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Update the confidential disclosure date (2021-05-13 -> 2021-05-14)
$disclosureText = $ws.Range("A80").Value2
$disclosureText = $disclosureText -replace [regex]::Escape("2021-05-13"), "2021-05-14"
$ws.Range("A80").Value2 = $disclosureText
$ws.Rows(80).AutoFit()

# Update Weight (D) and Percent Change (E) values per-row
$ws.Range("D2").Value = 0.06261980304155873
$ws.Range("E2").Value = 0.01984476274305846
$ws.Range("D3").Value = 0.03771774344353058
$ws.Range("E3").Value = 0.01943083439033111
$ws.Range("D4").Value = 0.03189401610007168
$ws.Range("E4").Value = 0.021067357939349
$ws.Range("D5").Value = 0.02936249500979318
$ws.Range("E5").Value = 0.01430231518727076
$ws.Range("D6").Value = 0.02659343876278042
$ws.Range("E6").Value = 0.02213508954527521
$ws.Range("D7").Value = 0.02569022455220805
$ws.Range("E7").Value = 0.01554179566563452
$ws.Range("D8").Value = 0.1907538941150909
$ws.Range("E8").Value = 0.004075761208343298
$ws.Range("D9").Value = 0.02500828032524015
$ws.Range("E9").Value = 0.001529771710990779
$ws.Range("D10").Value = 0.02300285208418945
$ws.Range("E10").Value = 0.00210572175428414
$ws.Range("D11").Value = 0.02203139572914094
$ws.Range("E11").Value = 0.01448385908923533
$ws.Range("D12").Value = 0.01985830496996005
$ws.Range("E12").Value = -0.02601771896377725
$ws.Range("D13").Value = 0.02029957235638622
$ws.Range("E13").Value = 0.01242829827915859
$ws.Range("D14").Value = 0.01732567037149663
$ws.Range("E14").Value = 0.007811011621261077
$ws.Range("D15").Value = 0.01591950837468313
$ws.Range("E15").Value = 0.02697922599598312
$ws.Range("D16").Value = 0.01439097978197621
$ws.Range("E16").Value = 0.03170181721713772
$ws.Range("D17").Value = 0.01429362735209416
$ws.Range("E17").Value = 0.0092592592592593
$ws.Range("D18").Value = 0.01455216009174416
$ws.Range("E18").Value = 0.001503040240486531
$ws.Range("D19").Value = 0.01335358783512033
$ws.Range("E19").Value = 0.03498656882657425
$ws.Range("D20").Value = 0.01320620707321555
$ws.Range("E20").Value = 0.02478920741989898
$ws.Range("D21").Value = 0.01268122773217268
$ws.Range("E21").Value = 0.0009313877677741278
$ws.Range("D22").Value = 0.01344720205895299
$ws.Range("E22").Value = -0.001457975986277726
$ws.Range("D23").Value = 0.01145386333533863
$ws.Range("E23").Value = 0.01080495528026226
$ws.Range("D24").Value = 0.01299814011523235
$ws.Range("E24").Value = 0.003035049931466638
$ws.Range("D25").Value = 0.01164721607802104
$ws.Range("E25").Value = -0.005500583860856834
$ws.Range("D26").Value = 0.008595169679439364
$ws.Range("E26").Value = 0.0325171654359373
$ws.Range("D27").Value = 0.009301452013877779
$ws.Range("E27").Value = 0.03122808817744938
$ws.Range("D28").Value = 0.01004042376717365
$ws.Range("E28").Value = 0.01722949689869058
$ws.Range("D29").Value = 0.009879243457405691
$ws.Range("E29").Value = 0.02429343735030454
$ws.Range("D30").Value = 0.009726493995323599
$ws.Range("E30").Value = 0.01856243355957132
$ws.Range("D31").Value = 0.008376723234446673
$ws.Range("E31").Value = 0.02481022032956881
$ws.Range("D32").Value = 0.01049855284897153
$ws.Range("E32").Value = 0.02247778358599062
$ws.Range("D33").Value = 0.009588935966335846
$ws.Range("E33").Value = -0.002040469307940884
$ws.Range("D34").Value = 0.009104599674042748
$ws.Range("E34").Value = 0.004035956705191523
$ws.Range("D35").Value = 0.009462831164393042
$ws.Range("E35").Value = 0.002311409960075617
$ws.Range("D36").Value = 0.008407384477681584
$ws.Range("E36").Value = 0.007322264793529154
$ws.Range("D37").Value = 0.008697294287256812
$ws.Range("E37").Value = 0.01618655692729765
$ws.Range("D38").Value = 0.006820516009714468
$ws.Range("E38").Value = 0.03157305532718779
$ws.Range("D39").Value = 0.008930335643101918
$ws.Range("E39").Value = -0.001995012468827717
$ws.Range("D40").Value = 0.008060685750675161
$ws.Range("E40").Value = 0.02615298087739015
$ws.Range("D41").Value = 0.006849666063269837
$ws.Range("E41").Value = 0.02559219693450987
$ws.Range("D42").Value = 0.007076901269293547
$ws.Range("E42").Value = 0.02391629297458908
$ws.Range("D43").Value = 0.008079456317221047
$ws.Range("E43").Value = 0.009598157153826392
$ws.Range("D44").Value = 0.007438791429394311
$ws.Range("E44").Value = 0.01736396976274213
$ws.Range("D45").Value = 0.007268345140801797
$ws.Range("E45").Value = 0.004705418891709723
$ws.Range("D46").Value = 0.008013441189114753
$ws.Range("E46").Value = 0.0157217722725107
$ws.Range("D47").Value = 0.007510453634724154
$ws.Range("E47").Value = -0.00508323802262034
$ws.Range("D48").Value = 0.007192944729422561
$ws.Range("E48").Value = 0.01271617497456767
$ws.Range("D49").Value = 0.006590974251021983
$ws.Range("E49").Value = 0.006425920897818793
$ws.Range("D50").Value = 0.007372776301288016
$ws.Range("E50").Value = 0.0007119971520115342
$ws.Range("D51").Value = 0.006682997748876178
$ws.Range("E51").Value = 0.008997375765401605
$ws.Range("D52").Value = 0.00667301594336131
$ws.Range("E52").Value = 0.02131133863335677
$ws.Range("D53").Value = 0.005245816595482606
$ws.Range("E53").Value = 0.05321810325221765
$ws.Range("D54").Value = 0.006201922444838932
$ws.Range("E54").Value = 0.006514825074381836
$ws.Range("D55").Value = 0.005665171731555414
$ws.Range("E55").Value = 0.01389912603980226
$ws.Range("D56").Value = 0.005701177814077967
$ws.Range("E56").Value = 0.01315635280035488
$ws.Range("D57").Value = 0.006861437435510476
$ws.Range("E57").Value = 0.00329206658320591
$ws.Range("D58").Value = 0.005543361889754466
$ws.Range("E58").Value = 0.007575757575757569
$ws.Range("D59").Value = 0.005389459151337987
$ws.Range("E59").Value = 0.01686810997476407
$ws.Range("D60").Value = 0.005044510222907874
$ws.Range("E60").Value = -0.003910191725529755
$ws.Range("D61").Value = 0.004835568365636511
$ws.Range("E61").Value = 0.03193414148724427
$ws.Range("D62").Value = 0.005074018189808399
$ws.Range("E62").Value = 0.01097264675915044
$ws.Range("D63").Value = 0.004255191992393594
$ws.Range("E63").Value = 0.003214953271028165
$ws.Range("D64").Value = 0.00413843270557427
$ws.Range("E64").Value = -0.002460024600246191
$ws.Range("D65").Value = 0.003882643968237115
$ws.Range("E65").Value = -0.006964929531301234
$ws.Range("D66").Value = 0.003844466544753958
$ws.Range("E66").Value = -0.002730883813306884
$ws.Range("D67").Value = 0.003868566043327701
$ws.Range("E67").Value = 0.01338432122370947
$ws.Range("D68").Value = 0.003619696963996868
$ws.Range("E68").Value = 0.01987475280158213
$ws.Range("D69").Value = 0.003626457549405344
$ws.Range("E69").Value = -0.00701831341155823
$ws.Range("D70").Value = 0.002922163622440012
$ws.Range("E70").Value = 0.03772455089820359
$ws.Range("D71").Value = 0.0028728511171076
$ws.Range("E71").Value = 0.02279900332225893
$ws.Range("D72").Value = 0.00217376681772174
$ws.Range("E72").Value = 0.04227877279961945
$ws.Range("D73").Value = 0.001935356761699314
$ws.Range("E73").Value = 0.01378786010767263
$ws.Range("D74").Value = 0.001885646574872286
$ws.Range("E74").Value = 0.02615151003880545
$ws.Range("D75").Value = 0.001398566280266335
$ws.Range("E75").Value = 0.08280254777070062
$ws.Range("D76").Value = 0.001671932539665528
$ws.Range("E76").Value = 0.03672517958232246
$ws.Range("E77").Value = 0.01171140584806696

$ws.Protect("D382")
